# Append the new daily price row (2024-11-07) to the end of the data table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row in column A (xlUp = -4162) and append right after it,
# matching the existing sheet's pattern of one row per day.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

$ws.Cells.Item($newRow, 1).Value = "2024-11-07 00:00:00"
$ws.Cells.Item($newRow, 2).Value = 74950
$ws.Cells.Item($newRow, 3).Value = 10407.41
$ws.Cells.Item($newRow, 4).Value = 9210.1
$ws.Cells.Item($newRow, 5).Value = 7.1618
